# Discharge_Aug05.xlsx - "lots of discharge data"
# Adds a new discharge-calc block (rows 34-47) to the "stn3" sheet, bolds the
# "new depth" label, and updates the active sheet/selection state so that
# "stn3" (not "stn1") is the tab that is focused when the workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("stn1")
$ws3 = $wb.Worksheets.Item("stn3")

# ---------------------------------------------------------------------------
# New block on stn3: label + header + 12 rows of discharge-segment data,
# mirroring the existing A19:F30 block but reading its B/C results back in.
# ---------------------------------------------------------------------------

$ws3.Range("A34").Value = "new depth"
$ws3.Range("A34").Font.Bold = $true

$ws3.Range("A35").Value = "X"
$ws3.Range("B35").Value = "V"
$ws3.Range("C35").Value = "D"
$ws3.Range("D35").Value = "segment"
$ws3.Range("E35").Value = "Q"
$ws3.Range("F35").Value = "Qtotal"

$ws3.Cells.Item(36,1).Value = 0.7
$ws3.Cells.Item(36,2).Value = 0
$ws3.Cells.Item(36,3).Formula = "=C19*2.54"
$ws3.Cells.Item(36,4).Formula = "=A36"
$ws3.Cells.Item(36,6).Formula = "=SUM(E36:E47)"

$ws3.Cells.Item(37,1).Value = 0.75
$ws3.Cells.Item(37,2).Value = 0.33748
$ws3.Cells.Item(37,3).Formula = "=C20*2.54"
$ws3.Cells.Item(37,4).Formula = "=(A37+(A38-A37)/2)"
$ws3.Cells.Item(37,5).Formula = "=(D37-D36)*(B37)*C37"

$ws3.Cells.Item(38,1).Value = 0.8
$ws3.Cells.Item(38,2).Value = 0.4004
$ws3.Cells.Item(38,3).Formula = "=C21*2.54"
$ws3.Cells.Item(38,4).Formula = "=(A38+(A39-A38)/2)"
$ws3.Cells.Item(38,5).Formula = "=(D38-D37)*(B38)*C38"

$ws3.Cells.Item(39,1).Value = 0.85
$ws3.Cells.Item(39,2).Value = 0.38324
$ws3.Cells.Item(39,3).Formula = "=C22*2.54"
$ws3.Cells.Item(39,4).Formula = "=(A39+(A40-A39)/2)"
$ws3.Cells.Item(39,5).Formula = "=(D39-D38)*(B39)*C39"

$ws3.Cells.Item(40,1).Value = 0.9
$ws3.Cells.Item(40,2).Value = 0.29172
$ws3.Cells.Item(40,3).Formula = "=C23*2.54"
$ws3.Cells.Item(40,4).Formula = "=(A40+(A41-A40)/2)"
$ws3.Cells.Item(40,5).Formula = "=(D40-D39)*(B40)*C40"

$ws3.Cells.Item(41,1).Value = 0.95
$ws3.Cells.Item(41,2).Value = 0.08008
$ws3.Cells.Item(41,3).Formula = "=C24*2.54"
$ws3.Cells.Item(41,4).Formula = "=(A41+(A42-A41)/2)"
$ws3.Cells.Item(41,5).Formula = "=(D41-D40)*(B41)*C41"

$ws3.Cells.Item(42,1).Value = 1
$ws3.Cells.Item(42,2).Value = 0.04004
$ws3.Cells.Item(42,3).Formula = "=C25*2.54"
$ws3.Cells.Item(42,4).Formula = "=(A42+(A43-A42)/2)"
$ws3.Cells.Item(42,5).Formula = "=(D42-D41)*(B42)*C42"

$ws3.Cells.Item(43,1).Value = 1.05
$ws3.Cells.Item(43,2).Value = 0.04576
$ws3.Cells.Item(43,3).Formula = "=C26*2.54"
$ws3.Cells.Item(43,4).Formula = "=(A43+(A44-A43)/2)"
$ws3.Cells.Item(43,5).Formula = "=(D43-D42)*(B43)*C43"

$ws3.Cells.Item(44,1).Value = 1.1
$ws3.Cells.Item(44,2).Value = 0.0858
$ws3.Cells.Item(44,3).Formula = "=C27*2.54"
$ws3.Cells.Item(44,4).Formula = "=(A44+(A45-A44)/2)"
$ws3.Cells.Item(44,5).Formula = "=(D44-D43)*(B44)*C44"

$ws3.Cells.Item(45,1).Value = 1.15
$ws3.Cells.Item(45,2).Value = 0.10868
$ws3.Cells.Item(45,3).Formula = "=C28*2.54"
$ws3.Cells.Item(45,4).Formula = "=(A45+(A46-A45)/2)"
$ws3.Cells.Item(45,5).Formula = "=(D45-D44)*(B45)*C45"

$ws3.Cells.Item(46,1).Value = 1.2
$ws3.Cells.Item(46,2).Value = 0.06863999999999999
$ws3.Cells.Item(46,3).Formula = "=C29*2.54"
$ws3.Cells.Item(46,4).Formula = "=(A46+(A47-A46)/2)"
$ws3.Cells.Item(46,5).Formula = "=(D46-D45)*(B46)*C46"

$ws3.Cells.Item(47,1).Value = 1.25
$ws3.Cells.Item(47,2).Value = 0
$ws3.Cells.Item(47,3).Formula = "=C30*2.54"
$ws3.Cells.Item(47,4).Formula = "=(A47+(A48-A47)/2)"
$ws3.Cells.Item(47,5).Formula = "=(D47-D46)*(B47)*C47"

# ---------------------------------------------------------------------------
# View/selection state: stn1 is no longer the selected tab, scrolled further
# down; stn3 becomes the selected tab with F36 (the new total) selected.
# ---------------------------------------------------------------------------

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

$ws3.Activate()
$ws3.Range("F36").Select()
